# <MS 3/1> Radio button & Auto Suggest code added in Utility
# This script rebuilds the "RedBus" sheet (Source/Destination lookup table)
# and nudges the "FBLogin" sheet's remembered selection, matching the
# captured OOXML diff.

$wb = $excel.ActiveWorkbook

$wsFBLogin = $wb.Worksheets.Item(4)
$wsRedBus  = $wb.Worksheets.Item(5)

# --- FBLogin: only the remembered selection moved (C10 -> C2) ---
$wsFBLogin.Range("C2").Select() | Out-Null

# --- RedBus: build the Source/Destination table ---

# Header row (reuse the existing bold/shaded header style by copying it
# from FBLogin's own header row, so no new style gets minted).
$wsFBLogin.Range("A1").Copy()
$wsRedBus.Range("A1").PasteSpecial(-4122)
$wsRedBus.Range("A1").Value = "TestCaseid"
$wsRedBus.Range("B1").Value = "Source"
$wsRedBus.Range("C1").Value = "Destination"

$wsFBLogin.Range("B1").Copy()
$wsRedBus.Range("B1").PasteSpecial(-4122)
$wsFBLogin.Range("C1").Copy()
$wsRedBus.Range("C1").PasteSpecial(-4122)

# Data rows: reuse FBLogin's data-row style (left/top/wrap, bordered) as a
# base, then flip vertical alignment to "center" to match the new style.
$wsFBLogin.Range("A2").Copy()
$wsRedBus.Range("A2:C4").PasteSpecial(-4122)

$wsRedBus.Range("A2:C4").VerticalAlignment = -4108

$wsRedBus.Range("A2").Value = "TC001"
$wsRedBus.Range("B2").Value = "Chennai"
$wsRedBus.Range("C2").Value = "Nagercoil"

$wsRedBus.Range("A3").Value = "TC002"
$wsRedBus.Range("B3").Value = "Chennai"
$wsRedBus.Range("C3").Value = "Nagercoil"

$wsRedBus.Range("A4").Value = "TC003"
$wsRedBus.Range("B4").Value = "Chennai"
$wsRedBus.Range("C4").Value = "Nagercoil"

# Row heights
$wsRedBus.Rows.Item(2).RowHeight = 34.5
$wsRedBus.Rows.Item(3).RowHeight = 31.5
$wsRedBus.Rows.Item(4).RowHeight = 29.25

# Column widths (closest values the host's pixel-quantized column-width
# model can reach to 15 / 21.140625 / 22.7109375 "characters")
$wsRedBus.Columns.Item(1).ColumnWidth = 14.15
$wsRedBus.Columns.Item(2).ColumnWidth = 20.33
$wsRedBus.Columns.Item(3).ColumnWidth = 21.83

# Page setup (portrait)
$wsRedBus.PageSetup.Orientation = 1

# Selection on RedBus moves to C7, and RedBus stays the active sheet/tab.
$wsRedBus.Range("C7").Select() | Out-Null
$wsRedBus.Activate() | Out-Null
